$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats/styles) from row 2 down to the new row 3
# so the new row matches the existing layout (date format, percentage format, etc.)
$ws.Range("A2:W2").Copy()
$ws.Range("A3:W3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new result row with the ticker's latest data
$ws.Cells.Item(3, 1).Value = 42632.878761574073
$ws.Cells.Item(3, 2).Value = -1
$ws.Cells.Item(3, 3).Value = "Neutral"
$ws.Cells.Item(3, 4).Value = 24
$ws.Cells.Item(3, 5).Value = 24523
$ws.Cells.Item(3, 6).Value = 2641
$ws.Cells.Item(3, 7).Value = 59
$ws.Cells.Item(3, 8).Value = 37
$ws.Cells.Item(3, 9).Value = 86
$ws.Cells.Item(3, 10).Value = 13
$ws.Cells.Item(3, 11).Value = 12086
$ws.Cells.Item(3, 12).Value = 373
$ws.Cells.Item(3, 13).Value = 237
$ws.Cells.Item(3, 14).Value = 40
$ws.Cells.Item(3, 15).Value = 6
$ws.Cells.Item(3, 16).Value = "Bag"
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = -31.57
$ws.Cells.Item(3, 19).Value = -0.0872
$ws.Cells.Item(3, 20).Value = -0.74
$ws.Cells.Item(3, 21).Value = 6.75
$ws.Cells.Item(3, 22).Value = 1.88
$ws.Cells.Item(3, 23).Value = 0
